# Removing one overlapping ID (keeping Kayna coded, deleting Jared coded)
# Update measureNAprops / measureNAsums values after removing one
# overlapping double-coded ID from the TrHO_C group (row 21), which
# shifts the proportions/sums for that row and a few related rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.2579
$ws.Range("B17").Value = 0.1171
$ws.Range("B18").Value = 0.1394
$ws.Range("B21").Value = 0.1399
$ws.Range("C21").Value = 282
$ws.Range("B25").Value = 0.0655
